$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.827.05"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "3.817.09"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.54"
$ws.Range("E5").Value = "  +1.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.19"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.32"
$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000255"
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").Value = "4.456.62"
$ws.Range("E14").Value = "  +0.90%  "

$ws.Range("D15").Value = "3.798.64"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.50"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "67.847.00"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.31"
$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("E21").Value = "  -2.33%  "

$ws.Range("E22").Value = "  +0.88%  "

$ws.Range("E23").Value = "  -3.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.45"
$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("E26").Value = "  -1.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").Value = "3.968.42"
$ws.Range("E29").Value = "  +0.91%  "

$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("E32").Value = "  +1.70%  "

$ws.Range("E33").Value = "  -0.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("D36").Value = "3.762.10"
$ws.Range("E36").Value = "  +0.66%  "

$ws.Range("E38").Value = "  +1.51%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.81"
$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.16"
$ws.Range("E44").Value = "  +2.17%  "

$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.52"
$ws.Range("E46").Value = "  +10.91%  "

$ws.Range("E47").Value = "  -5.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.39"
$ws.Range("E48").Value = "  +11.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.35"
$ws.Range("E49").Value = "  +0.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "148.67"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("E51").Value = "  +0.51%  "
